# Customize the Details page
# - Adds a new "Loeng" entry as the first log row (row 7), pushing the
#   previous first entry down to row 7 -> 8.
# - Adds a partial new entry (date + start time + activity) in row 9.
# - Updates the selected cell in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")
$ws.Activate()

# --- Push the existing first log entry (row 7) down into row 8 -----------
# Copy formatting + values from row 7 (B:I) down to row 8 (B:I) so the
# previously-empty row 8 picks up the exact same cell styles row 7 had.
$ws.Range("B7:I7").Copy($ws.Range("B8:I8"))

# --- Write the new first entry into row 7 ---------------------------------
$ws.Range("B7").Value2 = 43872
$ws.Range("C7").Value2 = 0.33333333333333331
$ws.Range("D7").Value2 = 0.39583333333333331
$ws.Range("F7").Value2 = 90
$ws.Range("G7").Value2 = "Loeng"

# --- Fill in the new (partial) third entry in row 9 -----------------------
$ws.Range("B9").Value2 = 43876
$ws.Range("C9").Value2 = 0.95833333333333337
$ws.Range("G9").Value2 = "Kodutöö 3"

# --- Update the active selection shown in the sheet view ------------------
$ws.Range("D9").Select() | Out-Null
